$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append two new weekly observations ---
$data = $wb.Worksheets.Item("Data")

# Clone the formatting of the last existing data row (border/number-format/
# alignment) onto the two new rows, then overwrite with the new values.
$data.Range("A109:B109").Copy()
$data.Range("A110:B110").PasteSpecial(-4122)
$data.Range("A109:B109").Copy()
$data.Range("A111:B111").PasteSpecial(-4122)

$data.Cells.Item(110, 1).Value = 45231
$data.Cells.Item(110, 2).Value = 820.487

$data.Cells.Item(111, 1).Value = 45238
$data.Cells.Item(111, 2).Value = 772.686

# --- Sheet "SeriesInfo": refresh the FRED series metadata ---
$info = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($cell, $text) {
    # Force plain text so date-looking strings (e.g. "2023-11-15") are not
    # auto-converted into date serials - matches the source data's literal
    # string cells - then drop the temporary "Text" number format so the
    # cell is left unstyled, same as before the edit.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $info.Cells.Item(3, 2) "2023-11-15"
Set-TextValue $info.Cells.Item(4, 2) "2023-11-15"
Set-TextValue $info.Cells.Item(7, 2) "2023-11-08"
Set-TextValue $info.Cells.Item(14, 2) "2023-11-09 15:38:01-06"
